$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '30.673.49'
$ws.Range("E2").Value = '  -0.12%  '

# Row 3
$ws.Range("D3").Value = '1.918.64'
$ws.Range("E3").Value = '  +1.36%  '

# Row 4
$ws.Range("E4").Value = '  +0.07%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.65'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.17%  '

# Row 6
$ws.Range("E6").Value = '  +0.09%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4933'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.27%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2974'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.61%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06759'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.50%  '

# Row 10
$ws.Range("D10").Value = '1.909.46'
$ws.Range("E10").Value = '  +0.91%  '

# Row 11
$ws.Range("E11").Value = '  -0.29%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07365'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.71%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.175'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.66%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.78'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.26%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6696'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.35%  '

# Row 16
$ws.Range("D16").Value = '30.650.81'
$ws.Range("E16").Value = '  -0.08%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000007938'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.61%  '

# Row 18
$ws.Range("E18").Value = '  +2.59%  '

# Row 19
$ws.Range("E19").Value = '  +0.05%  '

# Row 20
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.356'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +11.11%  '

# Row 21
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.117.71'
$ws.Range("E21").Value = '  -0.68%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.002'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.11%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '203.32'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +7.73%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.317'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.67%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.647'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.92%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.36'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.56%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.79'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.46%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.957'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.97%  '

# Row 29
$ws.Range("E29").Value = '  +6.29%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.371'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.86%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09181'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.20%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.065'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.30%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05270'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.18%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7413'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.20%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.117'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.78%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.724'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.80%  '

# Row 37
$ws.Range("E37").Value = '  +0.26%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.718'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.06%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9261'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.25%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.081'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.82%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4458'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.86%  '

# Row 42
$ws.Range("E42").Value = '  +26.25%  '

# Row 43
$ws.Range("E43").Value = '  +3.85%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '106.35'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.01%  '

# Row 45
$ws.Range("E45").Value = '  +0.17%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1392'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.70%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.639'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.61%  '

# Row 48
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.057'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.06%  '

# Row 49
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.32'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.11%  '

# Row 50
$ws.Range("E50").Value = '  +0.18%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4030'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.36%  '
